$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "L2" = 1.35
    "Q2" = 1.79
    "S2" = 3
    "T2" = 1.87
    "F4" = 2.5
    "G4" = 2.76
    "H4" = 2.74
    "I4" = 3.45
    "K4" = 3.8
    "O4" = 1.33
    "V4" = 1.44
    "W4" = 1.56
    "G5" = 2.98
    "J5" = 2.94
    "K5" = 4.6
    "M5" = 1.07
    "N5" = 2.48
    "O5" = 1.29
    "R5" = 1.23
    "G6" = 5.9
    "W6" = 1.2
    "S7" = 2.64
    "M8" = 1.06
    "H9" = 4.2
    "L9" = 1.45
    "AA10" = 450
    "AB10" = 9.2
    "AC10" = 13
    "AE10" = 180
    "AF10" = 8.4
    "AI10" = 150
    "AN10" = 5
    "F10" = 1.33
    "G10" = 1.35
    "H10" = 11.5
    "I10" = 12
    "J10" = 5.9
    "K10" = 6
    "Q10" = 1.65
    "S10" = 2.62
    "T10" = 2.02
    "U10" = 1.91
    "W10" = 3.85
    "X10" = 24
    "Y10" = 40
    "Z10" = 110
    "AG11" = 20
    "AO11" = 8
    "F11" = 5.4
    "G11" = 5.5
    "H11" = 1.7
    "I11" = 1.71
    "P11" = 2.32
    "Q11" = 1.72
    "U11" = 2.26
    "V11" = 2.4
    "X11" = 21
    "AC13" = 11
    "AM13" = 60
    "K13" = 3.8
    "T13" = 1.57
    "U13" = 2.54
    "V13" = 1.34
    "W13" = 1.85
    "AB14" = 11
    "AC14" = 7.6
    "AE14" = 38
    "AH14" = 17
    "AI14" = 46
    "AM14" = 85
    "AN14" = 18.5
    "AO14" = 34
    "F14" = 2.36
    "G14" = 2.38
    "I14" = 3.4
    "J14" = 3.55
    "L14" = 1.39
    "M14" = 1.07
    "N14" = 4.1
    "O14" = 1.31
    "P14" = 2
    "Q14" = 1.96
    "R14" = 1.39
    "S14" = 3.4
    "T14" = 1.75
    "U14" = 2.26
    "V14" = 1.42
    "X14" = 15.5
    "Y14" = 13.5
    "AA15" = 260
    "AE15" = 110
    "AM15" = 80
    "AN15" = 3.4
    "AO15" = 80
    "F15" = 1.34
    "G15" = 1.35
    "H15" = 9.6
    "I15" = 9.8
    "U15" = 2.44
    "X15" = 44
    "AE16" = 130
    "AI16" = 110
    "AK16" = 12.5
    "AL16" = 32
    "AN16" = 4.5
    "AO16" = 120
    "F16" = 1.37
    "G16" = 1.38
    "H16" = 10
    "I16" = 10.5
    "S16" = 2.32
    "T16" = 1.85
    "V16" = 1.1
    "W16" = 3.6
    "Z16" = 95
    "AA17" = 14.5
    "AB17" = 30
    "AG17" = 25
    "AH17" = 19.5
    "AJ17" = 190
    "AL17" = 70
    "AM17" = 85
    "AO17" = 5.3
    "F17" = 6.8
    "G17" = 7
    "H17" = 1.51
    "I17" = 1.52
    "J17" = 5.1
    "N17" = 5.8
    "P17" = 2.68
    "U17" = 2.28
    "V17" = 2.96
    "AH18" = 16
    "AK18" = 22
    "AO18" = 29
    "F18" = 2.3
    "G18" = 2.32
    "H18" = 3.35
    "N18" = 4.3
    "P18" = 2.12
    "Q18" = 1.86
    "S18" = 3.1
    "T18" = 1.71
    "V18" = 1.41
    "W18" = 1.75
    "AA19" = 40
    "AD19" = 12.5
    "G19" = 2.7
    "J19" = 3.7
    "L19" = 1.33
    "O19" = 1.23
    "P19" = 2.38
    "Q19" = 1.7
    "R19" = 1.56
    "S19" = 2.7
    "T19" = 1.59
    "V19" = 1.54
    "X19" = 20
    "AA20" = 120
    "AD20" = 19.5
    "AG20" = 9.6
    "AJ20" = 18
    "AL20" = 25
    "AM20" = 70
    "F20" = 1.71
    "G20" = 1.72
    "H20" = 5.2
    "I20" = 5.3
    "K20" = 4.4
    "O20" = 1.2
    "Q20" = 1.6
    "U20" = 2.48
    "V20" = 1.23
    "W20" = 2.38
    "X20" = 23
    "AC21" = 8.2
    "AN21" = 32
    "K21" = 3.75
    "N21" = 3.1
    "O21" = 1.33
    "P21" = 1.84
    "Q21" = 1.96
    "S21" = 3.1
    "U21" = 2.08
    "X21" = 16
    "F22" = 1.95
    "G22" = 2
    "I22" = 4.5
    "V22" = 1.29
    "W22" = 2
    "J23" = 1.03
    "K23" = 950
    "N23" = 1.28
    "O23" = 1.01
    "P23" = 1.28
    "R23" = 1.21
    "S23" = 1.28
    "T23" = 1.04
    "U23" = 1.04
    "AD24" = 14
    "AJ24" = 1000
    "AK24" = 1000
    "AO24" = 48
    "F24" = 2.6
    "G24" = 3
    "H24" = 2.74
    "I24" = 3.2
    "V24" = 1.45
    "W24" = 1.5
    "G25" = 1.9
    "Q25" = 2.14
    "W25" = 2.1
    "AA26" = 1000
    "AD26" = 24
    "AF26" = 11.5
    "AI26" = 1000
    "AJ26" = 29
    "AK26" = 60
    "AN26" = 34
    "F26" = 2.08
    "G26" = 2.2
    "I26" = 5.4
    "T26" = 2.34
    "V26" = 1.23
    "W26" = 1.83
    "Z26" = 40
    "W27" = 1.33
    "AN28" = 11
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Host "Applied $($updates.Count) cell updates"